$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the D-column "Example" links for the rows whose CodePen link was a
# placeholder ("[CodePen]()") with the real pen URLs (or N/A where there is
# no demo). Order matters for the order new strings land in the shared
# string table, so write them in the same order the source workbook does.
$ws.Range("D3").Value  = "[CodePen](https://codepen.io/maiCoding/pen/oQeQWe)"
$ws.Range("D10").Value = "N/A"
$ws.Range("D4").Value  = "[CodePen](https://codepen.io/maiCoding/pen/gQxZpJ)"
$ws.Range("D5").Value  = "[CodePen](https://codepen.io/maiCoding/pen/gQxZpJ)"
$ws.Range("D6").Value  = "[CodePen](https://codepen.io/maiCoding/pen/OajrdX)"

# Widen column D to fit the longer link text, and move the active selection.
$ws.Columns.Item(4).ColumnWidth = 51.1
$ws.Range("E14").Select()
